$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 142.9073533333333
$ws.Range("H2").Value = 428.72206
$ws.Range("I2").Value = 0.5576664151504187
$ws.Range("J2").Value = 0.5576664151504188
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 59.45197733333333
$ws.Range("N2").Value = 178.355932
$ws.Range("O2").Value = 0.304222453049858
$ws.Range("P2").Value = 0.304222453049858
$ws.Range("Q2").Value = 8496.124731139991
$ws.Range("R2").Value = 76465.12258025992
$ws.Range("S2").Value = 0.1696546448005808
$ws.Range("T2").Value = 0.1696546448005809

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 142.9073533333333
$ws.Range("H3").Value = 428.72206
$ws.Range("I3").Value = 0.5576664151504187
$ws.Range("J3").Value = 0.5576664151504188
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 107.1770123333333
$ws.Range("N3").Value = 321.531037
$ws.Range("O3").Value = 0.548436823552382
$ws.Range("P3").Value = 0.5484368235523819
$ws.Range("Q3").Value = 15316.38317073069
$ws.Range("R3").Value = 137847.4485365762
$ws.Range("S3").Value = 0.3058447973269396
$ws.Range("T3").Value = 0.3058447973269396

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 142.9073533333333
$ws.Range("H4").Value = 428.72206
$ws.Range("I4").Value = 0.5576664151504187
$ws.Range("J4").Value = 0.5576664151504188
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 28.793724
$ws.Range("N4").Value = 86.38117199999999
$ws.Range("O4").Value = 0.1473407233977601
$ws.Range("P4").Value = 0.1473407233977601
$ws.Range("Q4").Value = 4114.83488945048
$ws.Range("R4").Value = 37033.51400505432
$ws.Range("S4").Value = 0.08216697302289828
$ws.Range("T4").Value = 0.08216697302289828

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 63.967809
$ws.Range("H5").Value = 191.903427
$ws.Range("I5").Value = 0.2496211559306514
$ws.Range("J5").Value = 0.2496211559306514
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 59.45197733333333
$ws.Range("N5").Value = 178.355932
$ws.Range("O5").Value = 0.304222453049858
$ws.Range("P5").Value = 0.304222453049858
$ws.Range("Q5").Value = 3803.012730730996
$ws.Range("R5").Value = 34227.11457657896
$ws.Range("S5").Value = 0.07594036039036386
$ws.Range("T5").Value = 0.07594036039036386

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 63.967809
$ws.Range("H6").Value = 191.903427
$ws.Range("I6").Value = 0.2496211559306514
$ws.Range("J6").Value = 0.2496211559306514
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 107.1770123333333
$ws.Range("N6").Value = 321.531037
$ws.Range("O6").Value = 0.548436823552382
$ws.Range("P6").Value = 0.5484368235523819
$ws.Range("Q6").Value = 6855.87865412931
$ws.Range("R6").Value = 61702.90788716379
$ws.Range("S6").Value = 0.1369014338500803
$ws.Range("T6").Value = 0.1369014338500803

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 63.967809
$ws.Range("H7").Value = 191.903427
$ws.Range("I7").Value = 0.2496211559306514
$ws.Range("J7").Value = 0.2496211559306514
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 28.793724
$ws.Range("N7").Value = 86.38117199999999
$ws.Range("O7").Value = 0.1473407233977601
$ws.Range("P7").Value = 0.1473407233977601
$ws.Range("Q7").Value = 1841.871437230716
$ws.Range("R7").Value = 16576.84293507644
$ws.Range("S7").Value = 0.03677936169020724
$ws.Range("T7").Value = 0.03677936169020724

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 49.38440333333333
$ws.Range("H8").Value = 148.15321
$ws.Range("I8").Value = 0.1927124289189298
$ws.Range("J8").Value = 0.1927124289189298
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 59.45197733333333
$ws.Range("N8").Value = 178.355932
$ws.Range("O8").Value = 0.304222453049858
$ws.Range("P8").Value = 0.304222453049858
$ws.Range("Q8").Value = 2936.000427593524
$ws.Range("R8").Value = 26424.00384834172
$ws.Range("S8").Value = 0.05862744785891321
$ws.Range("T8").Value = 0.05862744785891322

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 49.38440333333333
$ws.Range("H9").Value = 148.15321
$ws.Range("I9").Value = 0.1927124289189298
$ws.Range("J9").Value = 0.1927124289189298
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 107.1770123333333
$ws.Range("N9").Value = 321.531037
$ws.Range("O9").Value = 0.548436823552382
$ws.Range("P9").Value = 0.5484368235523819
$ws.Range("Q9").Value = 5292.872805130974
$ws.Range("R9").Value = 47635.85524617877
$ws.Range("S9").Value = 0.1056905923753621
$ws.Range("T9").Value = 0.1056905923753621

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 49.38440333333333
$ws.Range("H10").Value = 148.15321
$ws.Range("I10").Value = 0.1927124289189298
$ws.Range("J10").Value = 0.1927124289189298
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 28.793724
$ws.Range("N10").Value = 86.38117199999999
$ws.Range("O10").Value = 0.1473407233977601
$ws.Range("P10").Value = 0.1473407233977601
$ws.Range("Q10").Value = 1421.96087948468
$ws.Range("R10").Value = 12797.64791536212
$ws.Range("S10").Value = 0.02839438868465454
$ws.Range("T10").Value = 0.02839438868465454
